$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,4
$arr[0,0] = 0.028184
$arr[0,1] = 8867.58
$arr[0,2] = 8853.5
$arr[0,3] = 0.026
$arr[1,0] = 0.17481
$arr[1,1] = 8867.57
$arr[1,2] = 8853.4
$arr[1,3] = 0.226
$arr[2,0] = 2
$arr[2,1] = 8867.469999999999
$arr[2,2] = 8851.700000000001
$arr[2,3] = 0.005
$arr[3,0] = 0.489477
$arr[3,1] = 8867
$arr[3,2] = 8851.6
$arr[3,3] = 0.045
$arr[4,0] = 1.229764
$arr[4,1] = 8866.82
$arr[4,2] = 8772.700000000001
$arr[4,3] = 0.426
$arr[5,0] = 0.131358
$arr[5,1] = 8866.66
$arr[5,2] = 8760.700000000001
$arr[5,3] = 1.58
$arr[6,0] = 2
$arr[6,1] = 8866.530000000001
$arr[6,2] = 8720.5
$arr[6,3] = 0.71
$arr[7,0] = 0.194997
$arr[7,1] = 8866.290000000001
$arr[7,2] = 8640.200000000001
$arr[7,3] = 0.3
$arr[8,0] = 0.053022
$arr[8,1] = 8866.129999999999
$arr[8,2] = 8640
$arr[8,3] = 0.064
$arr[9,0] = 2
$arr[9,1] = 8866.120000000001
$arr[9,2] = 8594.700000000001
$arr[9,3] = 0.43
$arr[10,0] = 1.557721
$arr[10,1] = 8865.200000000001
$arr[10,2] = 8451
$arr[10,3] = 0.191
$arr[11,0] = 1.25
$arr[11,1] = 8865.190000000001
$arr[11,2] = 8375.700000000001
$arr[11,3] = 0.28
$arr[12,0] = 0.240177
$arr[12,1] = 8865.040000000001
$arr[12,2] = 8375.4
$arr[12,3] = 0.008999999999999999
$arr[13,0] = 0.3
$arr[13,1] = 8865.030000000001
$arr[13,2] = 8360
$arr[13,3] = 0.475
$arr[14,0] = 0.0102
$arr[14,1] = 8864.440000000001
$arr[14,2] = 8300
$arr[14,3] = 0.09
$arr[15,0] = 0.141378
$arr[15,1] = 8864.18
$arr[15,2] = 8214.700000000001
$arr[15,3] = 0.013
$arr[16,0] = 0.0101
$arr[16,1] = 8862.620000000001
$arr[16,2] = 8211.1
$arr[16,3] = 0.3
$arr[17,0] = 0.183164
$arr[17,1] = 8861.799999999999
$arr[17,2] = 8211
$arr[17,3] = 0.9
$arr[18,0] = 0.499192
$arr[18,1] = 8861.49
$arr[18,2] = 8201
$arr[18,3] = 0.9
$arr[19,0] = 0.08876000000000001
$arr[19,1] = 8860.65
$arr[19,2] = 8199
$arr[19,3] = 0.837
$arr[20,0] = 2.4
$arr[20,1] = 8860.639999999999
$arr[20,2] = 8186.1
$arr[20,3] = 0.065
$arr[21,0] = 0.009599999999999999
$arr[21,1] = 8860.620000000001
$arr[21,2] = 8001
$arr[21,3] = 0.1
$arr[22,0] = 0.5
$arr[22,1] = 8860.530000000001
$arr[22,2] = 8000.4
$arr[22,3] = 0.027
$arr[23,0] = 0.054472
$arr[23,1] = 8860.459999999999
$arr[23,2] = 8000
$arr[23,3] = 2
$arr[24,0] = 0.226016
$arr[24,1] = 8860.379999999999
$arr[24,2] = 7860.5
$arr[24,3] = 0.15
$arr[25,0] = 0.670504
$arr[25,1] = 8860
$arr[25,2] = 7510.4
$arr[25,3] = 0.381
$arr[26,0] = 0.24
$arr[26,1] = 8859.889999999999
$arr[26,2] = 7502
$arr[26,3] = 0.1
$arr[27,0] = 2.4
$arr[27,1] = 8859.84
$arr[27,2] = 7501
$arr[27,3] = 0.079
$arr[28,0] = 0.683736
$arr[28,1] = 8859.709999999999
$arr[28,2] = 7403.6
$arr[28,3] = 0.15
$arr[29,0] = 0.152031
$arr[29,1] = 8858.709999999999
$arr[29,2] = 7400
$arr[29,3] = 0.066

$ws.Range("A2:D31").Value = $arr
